$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 3.034382333333333
$ws.Cells.Item(2,8).Value = 9.103147
$ws.Cells.Item(2,9).Value = 0.04287206612663805
$ws.Cells.Item(2,10).Value = 0.04287206612663806
$ws.Cells.Item(2,13).Value = 8.033114333333334
$ws.Cells.Item(2,14).Value = 24.099343
$ws.Cells.Item(2,15).Value = 0.1374088679258946
$ws.Cells.Item(2,16).Value = 0.1374088679258946
$ws.Cells.Item(2,17).Value = 24.37554021471345
$ws.Cells.Item(2,18).Value = 219.379861932421
$ws.Cells.Item(2,19).Value = 0.005891002072105427
$ws.Cells.Item(2,20).Value = 0.005891002072105429
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 3.034382333333333
$ws.Cells.Item(3,8).Value = 9.103147
$ws.Cells.Item(3,9).Value = 0.04287206612663805
$ws.Cells.Item(3,10).Value = 0.04287206612663806
$ws.Cells.Item(3,15).Value = 0.6355200716780686
$ws.Cells.Item(3,16).Value = 0.6355200716780686
$ws.Cells.Item(3,17).Value = 112.7375932738257
$ws.Cells.Item(3,18).Value = 1014.638339464431
$ws.Cells.Item(3,19).Value = 0.02724605853778791
$ws.Cells.Item(3,20).Value = 0.02724605853778791
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 3.034382333333333
$ws.Cells.Item(4,8).Value = 9.103147
$ws.Cells.Item(4,9).Value = 0.04287206612663805
$ws.Cells.Item(4,10).Value = 0.04287206612663806
$ws.Cells.Item(4,13).Value = 13.27489133333333
$ws.Cells.Item(4,14).Value = 39.824674
$ws.Cells.Item(4,15).Value = 0.2270710603960369
$ws.Cells.Item(4,16).Value = 0.2270710603960369
$ws.Cells.Item(4,17).Value = 40.28109573878645
$ws.Cells.Item(4,18).Value = 362.529861649078
$ws.Cells.Item(4,19).Value = 0.009735005516744716
$ws.Cells.Item(4,20).Value = 0.009735005516744718
$ws.Cells.Item(5,9).Value = 0.4108678777210459
$ws.Cells.Item(5,10).Value = 0.4108678777210459
$ws.Cells.Item(5,13).Value = 8.033114333333334
$ws.Cells.Item(5,14).Value = 24.099343
$ws.Cells.Item(5,15).Value = 0.1374088679258946
$ws.Cells.Item(5,16).Value = 0.1374088679258946
$ws.Cells.Item(5,17).Value = 233.6049409594594
$ws.Cells.Item(5,18).Value = 2102.444468635134
$ws.Cells.Item(5,19).Value = 0.0564568899447638
$ws.Cells.Item(5,20).Value = 0.05645688994476381
$ws.Cells.Item(6,9).Value = 0.4108678777210459
$ws.Cells.Item(6,10).Value = 0.4108678777210459
$ws.Cells.Item(6,15).Value = 0.6355200716780686
$ws.Cells.Item(6,16).Value = 0.6355200716780686
$ws.Cells.Item(6,19).Value = 0.261114783099495
$ws.Cells.Item(6,20).Value = 0.261114783099495
$ws.Cells.Item(7,9).Value = 0.4108678777210459
$ws.Cells.Item(7,10).Value = 0.4108678777210459
$ws.Cells.Item(7,13).Value = 13.27489133333333
$ws.Cells.Item(7,14).Value = 39.824674
$ws.Cells.Item(7,15).Value = 0.2270710603960369
$ws.Cells.Item(7,16).Value = 0.2270710603960369
$ws.Cells.Item(7,17).Value = 386.0371055966014
$ws.Cells.Item(7,18).Value = 3474.333950369412
$ws.Cells.Item(7,19).Value = 0.09329620467678709
$ws.Cells.Item(7,20).Value = 0.0932962046767871
$ws.Cells.Item(8,7).Value = 38.66298066666667
$ws.Cells.Item(8,8).Value = 115.988942
$ws.Cells.Item(8,9).Value = 0.546260056152316
$ws.Cells.Item(8,10).Value = 0.546260056152316
$ws.Cells.Item(8,13).Value = 8.033114333333334
$ws.Cells.Item(8,14).Value = 24.099343
$ws.Cells.Item(8,15).Value = 0.1374088679258946
$ws.Cells.Item(8,16).Value = 0.1374088679258946
$ws.Cells.Item(8,17).Value = 310.5841441627896
$ws.Cells.Item(8,18).Value = 2795.257297465106
$ws.Cells.Item(8,19).Value = 0.07506097590902534
$ws.Cells.Item(8,20).Value = 0.07506097590902536
$ws.Cells.Item(9,7).Value = 38.66298066666667
$ws.Cells.Item(9,8).Value = 115.988942
$ws.Cells.Item(9,9).Value = 0.546260056152316
$ws.Cells.Item(9,10).Value = 0.546260056152316
$ws.Cells.Item(9,15).Value = 0.6355200716780686
$ws.Cells.Item(9,16).Value = 0.6355200716780686
$ws.Cells.Item(9,17).Value = 1436.460837934107
$ws.Cells.Item(9,18).Value = 12928.14754140697
$ws.Cells.Item(9,19).Value = 0.3471592300407856
$ws.Cells.Item(9,20).Value = 0.3471592300407856
$ws.Cells.Item(10,7).Value = 38.66298066666667
$ws.Cells.Item(10,8).Value = 115.988942
$ws.Cells.Item(10,9).Value = 0.546260056152316
$ws.Cells.Item(10,10).Value = 0.546260056152316
$ws.Cells.Item(10,13).Value = 13.27489133333333
$ws.Cells.Item(10,14).Value = 39.824674
$ws.Cells.Item(10,15).Value = 0.2270710603960369
$ws.Cells.Item(10,16).Value = 0.2270710603960369
$ws.Cells.Item(10,17).Value = 513.2468669727676
$ws.Cells.Item(10,18).Value = 4619.221802754909
$ws.Cells.Item(10,19).Value = 0.124039850202505
$ws.Cells.Item(10,20).Value = 0.124039850202505
